$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores plain-looking decimals ("0.561", "1.00", ...) as
# TEXT in the source sheet (note trailing zeros like "1.00" / "0.860" that a
# real number would lose). Pre-setting NumberFormat to "@" (Text) before the
# assignment keeps Excel from silently re-interpreting these as numbers.

$ws.Range("D2").Value = '41.875.82'
$ws.Range("E2").Value = '  -0.49%  '
$ws.Range("D3").Value = '2.238.72'
$ws.Range("E3").Value = '  +0.22%  '
$ws.Range("E4").Value = '  -0.27%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '252.04'
$ws.Range("E5").Value = '  +9.06%  '
$ws.Range("E6").Value = '  +0.78%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '71.91'
$ws.Range("E7").Value = '  +2.33%  '
$ws.Range("E8").Value = '  -0.22%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.561'
$ws.Range("E9").Value = '  +1.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.14'
$ws.Range("E10").Value = '  +20.39%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0975'
$ws.Range("E11").Value = '  -0.59%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '58.48'
$ws.Range("E12").Value = '  +0.62%  '
$ws.Range("E13").Value = '  +0.59%  '
$ws.Range("E14").Value = '  +2.99%  '
$ws.Range("D15").Value = '2.574.11'
$ws.Range("E15").Value = '  +0.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.08'
$ws.Range("E16").Value = '  +3.13%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.860'
$ws.Range("E17").Value = '  -0.10%  '
$ws.Range("D18").Value = '2.241.47'
$ws.Range("E18").Value = '  -0.22%  '
$ws.Range("D19").Value = '41.826.59'
$ws.Range("E19").Value = '  -0.35%  '
$ws.Range("D20").Value = '0.0₃0970'
$ws.Range("E20").Value = '  -0.85%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.50'
$ws.Range("E23").Value = '  +0.78%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.26'
$ws.Range("E24").Value = '  +26.74%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.76'
$ws.Range("E26").Value = '  +3.72%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.50'
$ws.Range("E27").Value = '  +6.17%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.21'
$ws.Range("E28").Value = '  +2.82%  '
$ws.Range("E29").Value = '  +5.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '172.26'
$ws.Range("E30").Value = '  +3.92%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.81'
$ws.Range("E31").Value = '  +1.78%  '
$ws.Range("E32").Value = '  +1.87%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.126'
$ws.Range("E33").Value = '  +0.24%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.46'
$ws.Range("E34").Value = '  +4.11%  '
$ws.Range("E35").Value = '  +2.93%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.15'
$ws.Range("E38").Value = '  +18.02%  '
$ws.Range("E39").Value = '  +6.40%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.30'
$ws.Range("E40").Value = '  +4.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.10'
$ws.Range("E41").Value = '  +1.37%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '69.96'
$ws.Range("E42").Value = '  +9.68%  '
$ws.Range("E43").Value = '  +17.00%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.08'
$ws.Range("E44").Value = '  -0.39%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '11.64'
$ws.Range("E45").Value = '  +23.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.92'
$ws.Range("E46").Value = '  +13.77%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.93'
$ws.Range("E47").Value = '  +2.77%  '
$ws.Range("E48").Value = '  +3.23%  '
$ws.Range("E49").Value = '  -0.05%  '
$ws.Range("D50").Value = '0.0₃0153'
$ws.Range("E50").Value = '  +20.17%  '

$ws.Range("B21").Value = 'Litecoin'
$ws.Range("C21").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.51'
$ws.Range("E21").Value = '  +0.68%  '

$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.24'
$ws.Range("E22").Value = '  +1.46%  '

$ws.Range("B36").Value = 'InjectiveProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '26.67'
$ws.Range("E36").Value = '  +29.85%  '

$ws.Range("B37").Value = 'Filecoin'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.72'
$ws.Range("E37").Value = '  -0.79%  '

$ws.Range("B51").Value = 'ARBITRUM'
$ws.Range("C51").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.14'
$ws.Range("E51").Value = '  +5.79%  '

